$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster changes from MuSCs -> ECs, target cluster from MuSCs -> FAPs,
# and the numeric/statistical columns are recomputed for the new TPM data.
$ws.Cells.Item(2,1).Value2  = "ECs"
$ws.Cells.Item(2,2).Value2  = "Ccl21b"
$ws.Cells.Item(2,3).Value2  = "Ackr2"
$ws.Cells.Item(2,4).Value2  = "FAPs"
$ws.Cells.Item(2,5).Value2  = 1
$ws.Cells.Item(2,6).Value2  = 0.3333333333333333
$ws.Cells.Item(2,7).Value2  = 0.004739
$ws.Cells.Item(2,8).Value2  = 0.014217
$ws.Cells.Item(2,9).Value2  = 0.02588570741885795
$ws.Cells.Item(2,10).Value2 = 0.02588570741885795
$ws.Cells.Item(2,11).Value2 = 2
$ws.Cells.Item(2,12).Value2 = 0.6666666666666666
$ws.Cells.Item(2,13).Value2 = 0.2746273333333333
$ws.Cells.Item(2,14).Value2 = 0.823882
$ws.Cells.Item(2,15).Value2 = 1
$ws.Cells.Item(2,16).Value2 = 1
$ws.Cells.Item(2,17).Value2 = 0.001301458932666667
$ws.Cells.Item(2,18).Value2 = 0.011713130394
$ws.Cells.Item(2,19).Value2 = 0.02588570741885795
$ws.Cells.Item(2,20).Value2 = 0.02588570741885795

# Row 3: Sending cluster changes from MuSCs -> FAPs, target cluster from MuSCs -> FAPs,
# with refreshed numeric values.
$ws.Cells.Item(3,1).Value2  = "FAPs"
$ws.Cells.Item(3,2).Value2  = "Ccl21b"
$ws.Cells.Item(3,3).Value2  = "Ackr2"
$ws.Cells.Item(3,4).Value2  = "FAPs"
$ws.Cells.Item(3,5).Value2  = 2
$ws.Cells.Item(3,6).Value2  = 0.6666666666666666
$ws.Cells.Item(3,7).Value2  = 0.134289
$ws.Cells.Item(3,8).Value2  = 0.402867
$ws.Cells.Item(3,9).Value2  = 0.733523056250478
$ws.Cells.Item(3,10).Value2 = 0.733523056250478
$ws.Cells.Item(3,11).Value2 = 2
$ws.Cells.Item(3,12).Value2 = 0.6666666666666666
$ws.Cells.Item(3,13).Value2 = 0.2746273333333333
$ws.Cells.Item(3,14).Value2 = 0.823882
$ws.Cells.Item(3,15).Value2 = 1
$ws.Cells.Item(3,16).Value2 = 1
$ws.Cells.Item(3,17).Value2 = 0.03687942996600001
$ws.Cells.Item(3,18).Value2 = 0.331914869694
$ws.Cells.Item(3,19).Value2 = 0.733523056250478
$ws.Cells.Item(3,20).Value2 = 0.733523056250478

# Row 4: brand-new row for MuSCs -> FAPs (Ccl21b/Ackr2).
$ws.Cells.Item(4,1).Value2  = "MuSCs"
$ws.Cells.Item(4,2).Value2  = "Ccl21b"
$ws.Cells.Item(4,3).Value2  = "Ackr2"
$ws.Cells.Item(4,4).Value2  = "FAPs"
$ws.Cells.Item(4,5).Value2  = 2
$ws.Cells.Item(4,6).Value2  = 0.6666666666666666
$ws.Cells.Item(4,7).Value2  = 0.04404599999999999
$ws.Cells.Item(4,8).Value2  = 0.132138
$ws.Cells.Item(4,9).Value2  = 0.2405912363306641
$ws.Cells.Item(4,10).Value2 = 0.2405912363306641
$ws.Cells.Item(4,11).Value2 = 2
$ws.Cells.Item(4,12).Value2 = 0.6666666666666666
$ws.Cells.Item(4,13).Value2 = 0.2746273333333333
$ws.Cells.Item(4,14).Value2 = 0.823882
$ws.Cells.Item(4,15).Value2 = 1
$ws.Cells.Item(4,16).Value2 = 1
$ws.Cells.Item(4,17).Value2 = 0.012096235524
$ws.Cells.Item(4,18).Value2 = 0.108866119716
$ws.Cells.Item(4,19).Value2 = 0.2405912363306641
$ws.Cells.Item(4,20).Value2 = 0.2405912363306641
